$d = $word.ActiveDocument

$d.Content.Find.Execute("81×48=3888", $true, $false, $false, $false, $false, $true, 1, $false, "30×47=1410", 2)
$d.Content.Find.Execute("48×14=672", $true, $false, $false, $false, $false, $true, 1, $false, "32×62=1984", 2)
$d.Content.Find.Execute("98×27=2646", $true, $false, $false, $false, $false, $true, 1, $false, "75×39=2925", 2)
$d.Content.Find.Execute("21×72=1512", $true, $false, $false, $false, $false, $true, 1, $false, "44×89=3916", 2)
$d.Content.Find.Execute("53×67=3551", $true, $false, $false, $false, $false, $true, 1, $false, "51×62=3162", 2)
$d.Content.Find.Execute("28×70=1960", $true, $false, $false, $false, $false, $true, 1, $false, "11×15=165", 2)
$d.Content.Find.Execute("56×50=2800", $true, $false, $false, $false, $false, $true, 1, $false, "14×99=1386", 2)
$d.Content.Find.Execute("63×63=3969", $true, $false, $false, $false, $false, $true, 1, $false, "25×85=2125", 2)
$d.Content.Find.Execute("35×89=3115", $true, $false, $false, $false, $false, $true, 1, $false, "11×50=550", 2)
$d.Content.Find.Execute("27×71=1917", $true, $false, $false, $false, $false, $true, 1, $false, "17×68=1156", 2)
$d.Content.Find.Execute("18×41=738", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=1452", 2)
$d.Content.Find.Execute("54×38=2052", $true, $false, $false, $false, $false, $true, 1, $false, "53×96=5088", 2)
$d.Content.Find.Execute("26×36=936", $true, $false, $false, $false, $false, $true, 1, $false, "47×41=1927", 2)
$d.Content.Find.Execute("79×24=1896", $true, $false, $false, $false, $false, $true, 1, $false, "14×62=868", 2)
$d.Content.Find.Execute("48×55=2640", $true, $false, $false, $false, $false, $true, 1, $false, "88×56=4928", 2)
$d.Content.Find.Execute("55×94=5170", $true, $false, $false, $false, $false, $true, 1, $false, "56×99=5544", 2)
$d.Content.Find.Execute("75×64=4800", $true, $false, $false, $false, $false, $true, 1, $false, "72×16=1152", 2)
$d.Content.Find.Execute("57×25=1425", $true, $false, $false, $false, $false, $true, 1, $false, "49×52=2548", 2)
$d.Content.Find.Execute("89×86=7654", $true, $false, $false, $false, $false, $true, 1, $false, "11×54=594", 2)
$d.Content.Find.Execute("31×47=1457", $true, $false, $false, $false, $false, $true, 1, $false, "34×44=1496", 2)
$d.Content.Find.Execute("68×19=1292", $true, $false, $false, $false, $false, $true, 1, $false, "37×63=2331", 2)
$d.Content.Find.Execute("74×88=6512", $true, $false, $false, $false, $false, $true, 1, $false, "34×63=2142", 2)
$d.Content.Find.Execute("33×38=1254", $true, $false, $false, $false, $false, $true, 1, $false, "42×26=1092", 2)
$d.Content.Find.Execute("17×47=799", $true, $false, $false, $false, $false, $true, 1, $false, "49×71=3479", 2)
$d.Content.Find.Execute("94×84=7896", $true, $false, $false, $false, $false, $true, 1, $false, "81×95=7695", 2)
